$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "69.376.67"
Set-TextValue $ws "E2" "  -2.16%  "
Set-TextValue $ws "D3" "3.481.26"
Set-TextValue $ws "E3" "  -2.34%  "
Set-TextValue $ws "E4" "  -0.05%  "
Set-TextValue $ws "D5" "609.22"
Set-TextValue $ws "E5" "  +4.46%  "
Set-TextValue $ws "D6" "185.52"
Set-TextValue $ws "E6" "  -0.63%  "
Set-TextValue $ws "D7" "0.627"
Set-TextValue $ws "E7" "  -0.69%  "
Set-TextValue $ws "D9" "0.214"
Set-TextValue $ws "E9" "  -3.05%  "
Set-TextValue $ws "D10" "0.650"
Set-TextValue $ws "E10" "  -0.39%  "
Set-TextValue $ws "D11" "52.95"
Set-TextValue $ws "E11" "  -2.88%  "
Set-TextValue $ws "D12" "0.0000307"
Set-TextValue $ws "E12" "  -3.18%  "
Set-TextValue $ws "D13" "9.50"
Set-TextValue $ws "E13" "  +0.29%  "
Set-TextValue $ws "D14" "4.028.44"
Set-TextValue $ws "E14" "  -2.45%  "
Set-TextValue $ws "D15" "600.53"
Set-TextValue $ws "E15" "  +5.33%  "
Set-TextValue $ws "D16" "69.378.13"
Set-TextValue $ws "E16" "  -2.17%  "
Set-TextValue $ws "D17" "12.60"
Set-TextValue $ws "E17" "  +1.10%  "
Set-TextValue $ws "D18" "18.83"
Set-TextValue $ws "E18" "  -2.01%  "
Set-TextValue $ws "D19" "3.484.71"
Set-TextValue $ws "E19" "  -2.23%  "
Set-TextValue $ws "D21" "0.985"
Set-TextValue $ws "E21" "  -1.76%  "
Set-TextValue $ws "D22" "17.17"
Set-TextValue $ws "E22" "  -2.68%  "
Set-TextValue $ws "D23" "105.95"
Set-TextValue $ws "E23" "  +11.42%  "
Set-TextValue $ws "D24" "4.63"
Set-TextValue $ws "E24" "  +1.58%  "
Set-TextValue $ws "D25" "5.04"
Set-TextValue $ws "E25" "  +2.02%  "
Set-TextValue $ws "D26" "3.01"
Set-TextValue $ws "E26" "  +1.96%  "
Set-TextValue $ws "D27" "10.90"
Set-TextValue $ws "E27" "  -3.26%  "
Set-TextValue $ws "D28" "9.71"
Set-TextValue $ws "E28" "  +5.88%  "
Set-TextValue $ws "E29" "  +3.10%  "
Set-TextValue $ws "D30" "6.96"
Set-TextValue $ws "E30" "  -3.78%  "
Set-TextValue $ws "D31" "12.39"
Set-TextValue $ws "E31" "  +0.80%  "
Set-TextValue $ws "D32" "3.95"
Set-TextValue $ws "E32" "  +16.78%  "
Set-TextValue $ws "D33" "0.115"
Set-TextValue $ws "E33" "  -1.71%  "
Set-TextValue $ws "D34" "63.12"
Set-TextValue $ws "E34" "  -0.05%  "
Set-TextValue $ws "E35" "  -6.87%  "
Set-TextValue $ws "E36" "  -0.06%  "
Set-TextValue $ws "D37" "523.77"
Set-TextValue $ws "E37" "  -5.00%  "
Set-TextValue $ws "E38" "  -4.33%  "
Set-TextValue $ws "D39" "3.595.33"
Set-TextValue $ws "E39" "  +0.79%  "
Set-TextValue $ws "E40" "  +4.32%  "
Set-TextValue $ws "D41" "36.67"
Set-TextValue $ws "E41" "  -3.27%  "
Set-TextValue $ws "D42" "0.0₃0774"
Set-TextValue $ws "E42" "  -3.46%  "
Set-TextValue $ws "E43" "  +0.71%  "
Set-TextValue $ws "D44" "0.0461"
Set-TextValue $ws "E44" "  -1.03%  "
Set-TextValue $ws "D45" "2.93"
Set-TextValue $ws "E45" "  +0.74%  "
Set-TextValue $ws "E46" "  +2.83%  "
Set-TextValue $ws "E47" "  -4.54%  "
Set-TextValue $ws "D48" "8.80"
Set-TextValue $ws "E48" "  -5.69%  "
Set-TextValue $ws "E49" "  +0.42%  "
Set-TextValue $ws "D50" "0.000242"
Set-TextValue $ws "E50" "  -7.96%  "
Set-TextValue $ws "B51" "OceanProtocol"
Set-TextValue $ws "C51" "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
Set-TextValue $ws "D51" "1.35"
Set-TextValue $ws "E51" "  -10.74%  "
